$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 373.14285
$ws.Range("J38").Value = 1000
$ws.Range("L38").Value = 3000
$ws.Range("N38").Value = -3744
$ws.Range("H43").Value = 961.8570999999999
$ws.Range("I43").Value = 1244.6666
$ws.Range("J43").Value = 749.75
$ws.Range("K43").Value = 1244.6666
$ws.Range("L43").Value = 749.75
$ws.Range("M43").Value = -1175.6666
$ws.Range("N43").Value = -887.75
$ws.Range("L48").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("N48").Value = 0
$ws.Range("L56").ClearContents()
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("N56").Value = 0
$ws.Range("H58").Value = 472.42856
$ws.Range("J58").Value = 958.5
$ws.Range("L58").Value = 2875.5
$ws.Range("N58").Value = -3175.5
$ws.Range("H107").Value = 2777.8572
$ws.Range("I107").Value = 2282.5
$ws.Range("J107").Value = 5750
$ws.Range("K107").Value = 2282.5
$ws.Range("L107").Value = 5750
$ws.Range("M107").Value = -362.5
$ws.Range("N107").Value = -9590
$ws.Range("H112").Value = 2228.0256
$ws.Range("J112").Value = 2408.6177
$ws.Range("L112").Value = 7225.853099999999
$ws.Range("N112").Value = -9441.8531
$ws.Range("H137").Value = 1318.6177
$ws.Range("I137").Value = 978.0769
$ws.Range("J137").Value = 1529.4286
$ws.Range("K137").Value = 2934.2307
$ws.Range("L137").Value = 4588.2858
$ws.Range("M137").Value = -384.2307000000001
$ws.Range("N137").Value = -9688.2858

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4529.58
$ws.Range("I32").Value = 4501.8687
$ws.Range("J32").Value = 4617.3335
$ws.Range("K32").Value = 4501.8687
$ws.Range("L32").Value = 4617.3335
$ws.Range("M32").Value = -4214.8687
$ws.Range("N32").Value = -5191.3335
$ws.Range("H61").Value = 1427.56
$ws.Range("I61").Value = 1393.6818
$ws.Range("J61").Value = 1676
$ws.Range("K61").Value = 1393.6818
$ws.Range("L61").Value = 1676
$ws.Range("M61").Value = -1181.6818
$ws.Range("N61").Value = -2100
$ws.Range("H74").Value = 1469.1
$ws.Range("I74").Value = 739.5714
$ws.Range("K74").Value = 739.5714
$ws.Range("M74").Value = 134.4286
$ws.Range("H77").Value = 1469.1
$ws.Range("I77").Value = 739.5714
$ws.Range("K77").Value = 3697.857
$ws.Range("M77").Value = 670.143
$ws.Range("L111").ClearContents()
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("N111").Value = 0
$ws.Range("H132").Value = 1450.7441
$ws.Range("I132").Value = 1190.7333
$ws.Range("K132").Value = 3572.199900000001
$ws.Range("M132").Value = -1042.199900000001
$ws.Range("H136").Value = 1427.56
$ws.Range("I136").Value = 1393.6818
$ws.Range("J136").Value = 1676
$ws.Range("K136").Value = 4181.0454
$ws.Range("L136").Value = 5028
$ws.Range("M136").Value = -1631.0454
$ws.Range("N136").Value = -10128

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1679.5
$ws.Range("I31").Value = 1394.4
$ws.Range("K31").Value = 1394.4
$ws.Range("M31").Value = -1099.4
$ws.Range("H34").Value = 1679.5
$ws.Range("I34").Value = 1394.4
$ws.Range("K34").Value = 1394.4
$ws.Range("M34").Value = -1192.4
$ws.Range("H58").Value = 799.74286
$ws.Range("I58").Value = 715.34375
$ws.Range("K58").Value = 715.34375
$ws.Range("M58").Value = -512.34375
$ws.Range("H99").Value = 1812.2727
$ws.Range("I99").Value = 1705
$ws.Range("K99").Value = 1705
$ws.Range("M99").Value = -207
$ws.Range("H111").Value = 47999.668
$ws.Range("J111").Value = 47999.668
$ws.Range("L111").Value = 47999.668
$ws.Range("N111").Value = -56179.668
$ws.Range("H126").Value = 1812.2727
$ws.Range("I126").Value = 1705
$ws.Range("K126").Value = 5115
$ws.Range("M126").Value = -2645
$ws.Range("H134").Value = 1486.5862
$ws.Range("I134").Value = 1489.4584
$ws.Range("K134").Value = 4468.3752
$ws.Range("M134").Value = -1933.3752
$ws.Range("H136").Value = 799.74286
$ws.Range("I136").Value = 715.34375
$ws.Range("K136").Value = 2146.03125
$ws.Range("M136").Value = 403.96875

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 7740
$ws.Range("J101").Value = 7740
$ws.Range("L101").Value = 23220
$ws.Range("N101").Value = -28088
$ws.Range("H129").Value = 20834416
$ws.Range("I129").Value = 37037530
$ws.Range("J129").Value = 7577324.5
$ws.Range("K129").Value = 111112590
$ws.Range("L129").Value = 22731973.5
$ws.Range("M129").Value = -111107590
$ws.Range("N129").Value = -22741973.5
$ws.Range("H131").Value = 1170.2245
$ws.Range("J131").Value = 1176.8247
$ws.Range("L131").Value = 3530.474099999999
$ws.Range("N131").Value = -13610.4741
$ws.Range("H134").Value = 3434.7778
$ws.Range("I134").Value = 2603.375
$ws.Range("K134").Value = 7810.125
$ws.Range("M134").Value = -2740.125

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 21249.75
$ws.Range("I52").Value = 24999
$ws.Range("J52").Value = 20000
$ws.Range("K52").Value = 24999
$ws.Range("L52").Value = 20000
$ws.Range("M52").Value = -24740
$ws.Range("N52").Value = -20518
$ws.Range("H97").Value = 912
$ws.Range("I97").Value = 890
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 890
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -394
$ws.Range("N97").Value = -1992
$ws.Range("H122").Value = 1751.5714
$ws.Range("J122").Value = 4800
$ws.Range("L122").Value = 14400
$ws.Range("N122").Value = -19300
$ws.Range("H132").Value = 2063.8914
$ws.Range("I132").Value = 1656.0322
$ws.Range("K132").Value = 4968.096600000001
$ws.Range("M132").Value = -2438.096600000001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2100.5557
$ws.Range("I7").Value = 2050
$ws.Range("J7").Value = 2505
$ws.Range("K7").Value = 2050
$ws.Range("L7").Value = 2505
$ws.Range("M7").Value = -1938
$ws.Range("N7").Value = -2729
$ws.Range("H54").Value = 10084
$ws.Range("J54").Value = 10084
$ws.Range("L54").Value = 10084
$ws.Range("N54").Value = -11372
$ws.Range("H68").Value = 2188.5
$ws.Range("I68").Value = 1802
$ws.Range("J68").Value = 2729.6
$ws.Range("K68").Value = 1802
$ws.Range("L68").Value = 2729.6
$ws.Range("M68").Value = -1053
$ws.Range("N68").Value = -4227.6
$ws.Range("H71").Value = 2188.5
$ws.Range("I71").Value = 1802
$ws.Range("J71").Value = 2729.6
$ws.Range("K71").Value = 9010
$ws.Range("L71").Value = 13648
$ws.Range("M71").Value = -5266
$ws.Range("N71").Value = -21136
$ws.Range("M93").ClearContents()
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("H122").Value = 25012450
$ws.Range("I122").Value = 62526750
$ws.Range("J122").Value = 2917.5
$ws.Range("K122").Value = 187580250
$ws.Range("L122").Value = 8752.5
$ws.Range("M122").Value = -187577800
$ws.Range("N122").Value = -13652.5
$ws.Range("H126").Value = 2100.5557
$ws.Range("I126").Value = 2050
$ws.Range("J126").Value = 2505
$ws.Range("K126").Value = 6150
$ws.Range("L126").Value = 7515
$ws.Range("M126").Value = -3680
$ws.Range("N126").Value = -12455
$ws.Range("H135").Value = 34417.4
$ws.Range("J135").Value = 34417.4
$ws.Range("L135").Value = 34417.4
$ws.Range("N135").Value = -44557.4
$ws.Range("H136").Value = 12868.667
$ws.Range("I136").Value = 18384.666
$ws.Range("J136").Value = 1836.6666
$ws.Range("K136").Value = 55153.99800000001
$ws.Range("L136").Value = 5509.9998
$ws.Range("M136").Value = -52603.99800000001
$ws.Range("N136").Value = -10609.9998

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 18572704
$ws.Range("I122").Value = 21668014
$ws.Range("J122").Value = 847.5
$ws.Range("K122").Value = 65004042
$ws.Range("L122").Value = 2542.5
$ws.Range("M122").Value = -65001592
$ws.Range("N122").Value = -7442.5
